$d = $word.ActiveDocument

# Remove the entire first paragraph (the "I personally examined the
# patient separately..." attestation line), including its paragraph
# mark, so the following "OBJECTIVE:" paragraph becomes the first
# paragraph in the document body.
$d.Paragraphs(1).Range.Delete()
